# Weekly update: a new Betarraga price record is added at the top of the
# data block (row 91), pushing every existing record down by one row.
# The record that falls off the bottom (old row 193) lands in the new
# last row (194).
#
# Net effect per the target diff:
#   - Insert a new row at row 91 (shifts rows 91-193 down to 92-194).
#   - The freshly inserted row 91 inherits all the "constant" columns
#     (Mercado, Region, Codreg, Categoria, Variedad, Calidad, Precio
#     minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg o Unidades,
#     Clasificacion) from what is now row 92 (i.e. the record that used
#     to be row 91).
#   - Only the Fecha (D) and Volumen (J) of the new row 91 are genuinely
#     new values: 2021-10-27 (serial 44494) and 4000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 91; existing rows 91-193 shift down to 92-194.
$ws.Rows.Item(91).Insert()

# Duplicate the row that is now directly below (row 92, the former row 91)
# into the newly blank row 91 so every column starts out identical.
$ws.Range("A92:R92").Copy()
$ws.Range("A91").PasteSpecial()

# Overwrite the two columns that actually hold new data for this week.
$ws.Range("D91").Value = 44494
$ws.Range("J91").Value = 4000
